# Auto-generated script applying scheduled-runner Leve market-data refresh
# to the Lamia_Profits workbook. Updates currentAveragePrice / Leve price /
# profit columns (H:N) for the affected leve rows across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9441.6
$ws.Range("I62").Value = 7874.364
$ws.Range("K62").Value = 7874.364
$ws.Range("M62").Value = -7250.364
$ws.Range("H65").Value = 9441.6
$ws.Range("I65").Value = 7874.364
$ws.Range("K65").Value = 39371.82
$ws.Range("M65").Value = -36251.82
$ws.Range("H76").Value = 7120
$ws.Range("I76").Value = 4198.6665
$ws.Range("K76").Value = 4198.6665
$ws.Range("M76").Value = -3883.6665
$ws.Range("H79").Value = 7120
$ws.Range("I79").Value = 4198.6665
$ws.Range("K79").Value = 4198.6665
$ws.Range("M79").Value = -3106.6665
$ws.Range("H138").Value = 2720.0312
$ws.Range("J138").Value = 3480.7778
$ws.Range("L138").Value = 10442.3334
$ws.Range("N138").Value = -20722.3334
$ws.Range("H141").Value = 5883.8696
$ws.Range("I141").Value = 3619.25
$ws.Range("K141").Value = 10857.75
$ws.Range("M141").Value = -5677.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 40002852
$ws.Range("I45").Value = 62501550
$ws.Range("J45").Value = 5161.3335
$ws.Range("K45").Value = 62501550
$ws.Range("L45").Value = 5161.3335
$ws.Range("M45").Value = -62501173
$ws.Range("N45").Value = -5915.3335
$ws.Range("H132").Value = 4697.7896
$ws.Range("I132").Value = 3543.4
$ws.Range("J132").Value = 6917.769
$ws.Range("K132").Value = 10630.2
$ws.Range("L132").Value = 20753.307
$ws.Range("M132").Value = -8100.200000000001
$ws.Range("N132").Value = -25813.307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3907.5833
$ws.Range("I20").Value = 3963
$ws.Range("J20").Value = 3889.111
$ws.Range("K20").Value = 3963
$ws.Range("L20").Value = 3889.111
$ws.Range("M20").Value = -3716
$ws.Range("N20").Value = -4383.111
$ws.Range("H22").Value = 354.4
$ws.Range("I22").Value = 393
$ws.Range("K22").Value = 393
$ws.Range("M22").Value = -220
$ws.Range("H86").Value = 4143.722
$ws.Range("I86").Value = 2011
$ws.Range("K86").Value = 2011
$ws.Range("M86").Value = -888
$ws.Range("H89").Value = 4143.722
$ws.Range("I89").Value = 2011
$ws.Range("K89").Value = 10055
$ws.Range("M89").Value = -4439
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("M94").Value = -2549
$ws.Range("H99").Value = 1405.3636
$ws.Range("I99").Value = 1405.3636
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1405.3636
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 92.63640000000009
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 20749.25
$ws.Range("H21").Value = 3999
$ws.Range("J21").Value = 3999
$ws.Range("L21").Value = 3999
$ws.Range("N21").Value = -4469
$ws.Range("H22").Value = 1283.3334
$ws.Range("I22").Value = 1178.7693
$ws.Range("K22").Value = 1178.7693
$ws.Range("M22").Value = -828.7692999999999
$ws.Range("H134").Value = 4071.7144
$ws.Range("J134").Value = 5751.5
$ws.Range("L134").Value = 17254.5
$ws.Range("N134").Value = -22324.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 4474.5
$ws.Range("I20").Value = 4949
$ws.Range("J20").Value = 4000
$ws.Range("K20").Value = 14847
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = -14620
$ws.Range("N20").Value = -12454
$ws.Range("H50").Value = 10460.429
$ws.Range("I50").Value = 1644.6
$ws.Range("K50").Value = 4933.799999999999
$ws.Range("M50").Value = -4452.799999999999
$ws.Range("H53").Value = 10460.429
$ws.Range("I53").Value = 1644.6
$ws.Range("K53").Value = 4933.799999999999
$ws.Range("M53").Value = -4452.799999999999
$ws.Range("H132").Value = 4431.5884
$ws.Range("I132").Value = 2871.889
$ws.Range("J132").Value = 6186.25
$ws.Range("K132").Value = 25847.001
$ws.Range("L132").Value = 55676.25
$ws.Range("M132").Value = -23317.001
$ws.Range("N132").Value = -60736.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 18390
$ws.Range("I57").Value = 1000
$ws.Range("J57").Value = 29983.334
$ws.Range("K57").Value = 1000
$ws.Range("L57").Value = 29983.334
$ws.Range("M57").Value = -180
$ws.Range("N57").Value = -31623.334
$ws.Range("H70").Value = 5513
$ws.Range("I70").Value = 5513
$ws.Range("K70").Value = 5513
$ws.Range("M70").Value = -5243
$ws.Range("H73").Value = 5513
$ws.Range("I73").Value = 5513
$ws.Range("K73").Value = 5513
$ws.Range("M73").Value = -4577
$ws.Range("H80").Value = 5030.857
$ws.Range("I80").Value = 4859.222
$ws.Range("J80").Value = 5339.8
$ws.Range("K80").Value = 4859.222
$ws.Range("L80").Value = 5339.8
$ws.Range("M80").Value = -3861.222
$ws.Range("N80").Value = -7335.8
$ws.Range("H83").Value = 5030.857
$ws.Range("I83").Value = 4859.222
$ws.Range("J83").Value = 5339.8
$ws.Range("K83").Value = 24296.11
$ws.Range("L83").Value = 26699
$ws.Range("M83").Value = -19304.11
$ws.Range("N83").Value = -36683
$ws.Range("H102").Value = 3592.125
$ws.Range("I102").Value = 881.6667
$ws.Range("K102").Value = 881.6667
$ws.Range("M102").Value = 740.3333
$ws.Range("H122").Value = 5281.364
$ws.Range("I122").Value = 4355.4443
$ws.Range("J122").Value = 9448
$ws.Range("K122").Value = 13066.3329
$ws.Range("L122").Value = 28344
$ws.Range("M122").Value = -10616.3329
$ws.Range("N122").Value = -33244
$ws.Range("H125").Value = 33989
$ws.Range("J125").Value = 33989
$ws.Range("L125").Value = 33989
$ws.Range("N125").Value = -38909
$ws.Range("H126").Value = 14507
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 14507
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 43521
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -48461
$ws.Range("H132").Value = 8381.200000000001
$ws.Range("J132").Value = 17007
$ws.Range("L132").Value = 51021
$ws.Range("N132").Value = -56081.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3254.3572
$ws.Range("I81").Value = 2051
$ws.Range("J81").Value = 7666.6665
$ws.Range("K81").Value = 4102
$ws.Range("L81").Value = 15333.333
$ws.Range("M81").Value = -3041
$ws.Range("N81").Value = -17455.333
$ws.Range("H84").Value = 3254.3572
$ws.Range("I84").Value = 2051
$ws.Range("J84").Value = 7666.6665
$ws.Range("K84").Value = 20510
$ws.Range("L84").Value = 76666.66500000001
$ws.Range("M84").Value = -15206
$ws.Range("N84").Value = -87274.66500000001
$ws.Range("H122").Value = 6320.077
$ws.Range("I122").Value = 1408.1765
$ws.Range("J122").Value = 15598.111
$ws.Range("K122").Value = 4224.529500000001
$ws.Range("L122").Value = 46794.333
$ws.Range("M122").Value = -1774.529500000001
$ws.Range("N122").Value = -51694.333
